# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# per-language handback sheets, reflecting a fresh report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 00:51:46"
$wsZhCn.Range("H2").Value = "2016-03-18 00:52:04"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 00:51:49"
$wsDeDe.Range("H2").Value = "2016-03-18 00:52:12"
